$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 74607.5
$ws.Range("J28").Value = 104140.6
$ws.Range("L28").Value = 104140.6
$ws.Range("N28").Value = -105110.6
$ws.Range("H54").Value = 484038
$ws.Range("J54").Value = 950000
$ws.Range("L54").Value = 950000
$ws.Range("N54").Value = -950972
$ws.Range("H86").Value = 5782.2354
$ws.Range("J86").Value = 6592.7856
$ws.Range("L86").Value = 6592.7856
$ws.Range("N86").Value = -8838.785599999999
$ws.Range("H89").Value = 5782.2354
$ws.Range("J89").Value = 6592.7856
$ws.Range("L89").Value = 32963.928
$ws.Range("N89").Value = -44195.928
$ws.Range("H113").Value = 6961.8335
$ws.Range("J113").Value = 7714.4
$ws.Range("L113").Value = 7714.4
$ws.Range("N113").Value = -14222.4
$ws.Range("H132").Value = 3151.4348
$ws.Range("I132").Value = 3067.5789
$ws.Range("J132").Value = 3549.75
$ws.Range("K132").Value = 9202.736699999999
$ws.Range("L132").Value = 10649.25
$ws.Range("M132").Value = -6672.736699999999
$ws.Range("N132").Value = -15709.25
$ws.Range("H135").Value = 2859.182
$ws.Range("I135").Value = 1004
$ws.Range("J135").Value = 3919.2856
$ws.Range("K135").Value = 9036
$ws.Range("L135").Value = 35273.5704
$ws.Range("M135").Value = -6501
$ws.Range("N135").Value = -40343.5704

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 29416270
$ws.Range("I61").Value = 1550.4166
$ws.Range("K61").Value = 1550.4166
$ws.Range("M61").Value = -1338.4166
$ws.Range("H62").Value = 43478.75
$ws.Range("J62").Value = 43478.75
$ws.Range("L62").Value = 43478.75
$ws.Range("N62").Value = -44726.75
$ws.Range("H65").Value = 43478.75
$ws.Range("J65").Value = 43478.75
$ws.Range("L65").Value = 130436.25
$ws.Range("N65").Value = -136676.25
$ws.Range("H74").Value = 5879.154
$ws.Range("I74").Value = 1792.3334
$ws.Range("K74").Value = 1792.3334
$ws.Range("M74").Value = -918.3334
$ws.Range("H77").Value = 5879.154
$ws.Range("I77").Value = 1792.3334
$ws.Range("K77").Value = 8961.666999999999
$ws.Range("M77").Value = -4593.666999999999
$ws.Range("H102").Value = 1569.1
$ws.Range("I102").Value = 1496.5555
$ws.Range("K102").Value = 1496.5555
$ws.Range("M102").Value = 125.4445000000001
$ws.Range("H122").Value = 7350.865
$ws.Range("I122").Value = 7329.3193
$ws.Range("K122").Value = 21987.9579
$ws.Range("M122").Value = -19537.9579
$ws.Range("H136").Value = 29416270
$ws.Range("I136").Value = 1550.4166
$ws.Range("K136").Value = 4651.2498
$ws.Range("M136").Value = -2101.2498

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2008.8
$ws.Range("I20").Value = 2066.6843
$ws.Range("J20").Value = 909
$ws.Range("K20").Value = 2066.6843
$ws.Range("L20").Value = 909
$ws.Range("M20").Value = -1819.6843
$ws.Range("N20").Value = -1403
$ws.Range("H35").Value = 70000
$ws.Range("J35").Value = 70000
$ws.Range("L35").Value = 70000
$ws.Range("N35").Value = -70620
$ws.Range("H86").Value = 1948.9048
$ws.Range("I86").Value = 1838.2632
$ws.Range("K86").Value = 1838.2632
$ws.Range("M86").Value = -715.2632000000001
$ws.Range("H89").Value = 1948.9048
$ws.Range("I89").Value = 1838.2632
$ws.Range("K89").Value = 9191.316000000001
$ws.Range("M89").Value = -3575.316000000001
$ws.Range("H132").Value = 77998
$ws.Range("J132").Value = 77998
$ws.Range("L132").Value = 77998
$ws.Range("N132").Value = -88118

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5717.0264
$ws.Range("I31").Value = 2376.4348
$ws.Range("J31").Value = 10839.267
$ws.Range("K31").Value = 2376.4348
$ws.Range("L31").Value = 10839.267
$ws.Range("M31").Value = -2081.4348
$ws.Range("N31").Value = -11429.267
$ws.Range("H34").Value = 5717.0264
$ws.Range("I34").Value = 2376.4348
$ws.Range("J34").Value = 10839.267
$ws.Range("K34").Value = 2376.4348
$ws.Range("L34").Value = 10839.267
$ws.Range("M34").Value = -2174.4348
$ws.Range("N34").Value = -11243.267
$ws.Range("H107").Value = 812.8387
$ws.Range("J107").Value = 1096
$ws.Range("L107").Value = 1096
$ws.Range("N107").Value = -4936
$ws.Range("H110").Value = 199500
$ws.Range("J110").Value = 199500
$ws.Range("L110").Value = 199500
$ws.Range("N110").Value = -207680
$ws.Range("H132").Value = 5050.5864
$ws.Range("J132").Value = 15097.4
$ws.Range("L132").Value = 45292.2
$ws.Range("N132").Value = -50352.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 4269.1113
$ws.Range("I137").Value = 2711
$ws.Range("J137").Value = 5515.6
$ws.Range("K137").Value = 8133
$ws.Range("L137").Value = 16546.8
$ws.Range("M137").Value = -3033
$ws.Range("N137").Value = -26746.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 32500
$ws.Range("J57").Value = 55000
$ws.Range("L57").Value = 55000
$ws.Range("N57").Value = -56640
$ws.Range("H97").Value = 1495.3572
$ws.Range("I97").Value = 1415.4166
$ws.Range("J97").Value = 1975
$ws.Range("K97").Value = 1415.4166
$ws.Range("L97").Value = 1975
$ws.Range("M97").Value = -919.4166
$ws.Range("N97").Value = -2967
$ws.Range("H122").Value = 5563.778
$ws.Range("I122").Value = 4639.44
$ws.Range("J122").Value = 7664.5454
$ws.Range("K122").Value = 13918.32
$ws.Range("L122").Value = 22993.6362
$ws.Range("M122").Value = -11468.32
$ws.Range("N122").Value = -27893.6362
$ws.Range("H126").Value = 4678.7896
$ws.Range("I126").Value = 2223
$ws.Range("K126").Value = 6669
$ws.Range("M126").Value = -4199
$ws.Range("H132").Value = 3552.2703
$ws.Range("I132").Value = 2451.7188
$ws.Range("K132").Value = 7355.1564
$ws.Range("M132").Value = -4825.1564
$ws.Range("H136").Value = 22715.215
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 22715.215
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 68145.645
$ws.Range("M136").Value = $null
$ws.Range("N136").Value = -73245.645

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3293.1052
$ws.Range("I40").Value = 2454.375
$ws.Range("J40").Value = 7766.3335
$ws.Range("K40").Value = 2454.375
$ws.Range("L40").Value = 7766.3335
$ws.Range("M40").Value = -2318.375
$ws.Range("N40").Value = -8038.3335
$ws.Range("H100").Value = 4199.6
$ws.Range("I100").Value = 4199.6
$ws.Range("K100").Value = 4199.6
$ws.Range("M100").Value = -3658.6
$ws.Range("H132").Value = 10605.529
$ws.Range("J132").Value = 12874.75
$ws.Range("L132").Value = 38624.25
$ws.Range("N132").Value = -43684.25
$ws.Range("H136").Value = 16135005
$ws.Range("I136").Value = 3878.7083
$ws.Range("J136").Value = 26323084
$ws.Range("K136").Value = 11636.1249
$ws.Range("L136").Value = 78969252
$ws.Range("M136").Value = -9086.124899999999
$ws.Range("N136").Value = -78974352

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").Value = $null
$ws.Range("H45").Value = 8120.6665
$ws.Range("I45").Value = 6249.5
$ws.Range("J45").Value = 9056.25
$ws.Range("K45").Value = 6249.5
$ws.Range("L45").Value = 9056.25
$ws.Range("M45").Value = -5758.5
$ws.Range("N45").Value = -10038.25
$ws.Range("H100").Value = 407.41666
$ws.Range("I100").Value = 367.375
$ws.Range("K100").Value = 734.75
$ws.Range("M100").Value = -193.75
$ws.Range("H122").Value = 2445.75
$ws.Range("I122").Value = 2337.6316
$ws.Range("J122").Value = 4500
$ws.Range("K122").Value = 7012.8948
$ws.Range("L122").Value = 13500
$ws.Range("M122").Value = -4562.8948
$ws.Range("N122").Value = -18400

